$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns retain their text formatting for the affected rows
# so numeric-looking strings (e.g. "1.000", "27.700.01") are preserved exactly as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.700.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.903.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5191'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3774'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07240'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9002'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07629'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.879.51'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.445'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008695'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.31%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '27.740.45'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.141'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.122.63'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.579'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.882'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.162'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.41'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.855'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08976'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.851'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.36%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.229'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7703'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.648'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02083'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.065'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.65%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5509'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05284'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.657'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '114.34'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.515'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1509'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4797'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.47'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9989'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.612'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.66'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05988'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.06%  '
